# Update daily COVID-19 figures for Valais (rows 456-461).
# Only the manually-entered input columns are touched:
#   C = Nb nouveaux cas positifs
#   E = Patients COVID-19 intubes
#   F = Patients COVID-19 hospitalises hors SI
#   G = Patients COVID-19 aux SI total (y.c. intubes)
#   L = Nb nouveaux deces a l'hopital
#   M = Nb nouveaux deces extra-hospitaliers
# Columns B, H, J, K are driven by existing shared formulas and
# recalculate automatically once their precedents change.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 456 (2021-05-26): new positive cases revised 43 -> 44
$ws.Range("C456").Value = 44

# Row 457 (2021-05-27): new positive cases revised 52 -> 56
$ws.Range("C457").Value = 56

# Row 458 (2021-05-28): fill in hospital occupancy that was missing
$ws.Range("C458").Value = 26
$ws.Range("E458").Value = 3
$ws.Range("F458").Value = 3
$ws.Range("G458").Value = 12

# Row 459 (2021-05-29): newly reported day
$ws.Range("C459").Value = 20
$ws.Range("E459").Value = 4
$ws.Range("F459").Value = 3
$ws.Range("G459").Value = 12

# Row 460 (2021-05-30): newly reported day
$ws.Range("C460").Value = 10
$ws.Range("E460").Value = 4
$ws.Range("F460").Value = 3
$ws.Range("G460").Value = 13

# Row 461 (2021-05-31): newly reported day
$ws.Range("C461").Value = 2
$ws.Range("E461").Value = 4
$ws.Range("F461").Value = 2
$ws.Range("G461").Value = 11

# Columns L and M (rows 459-461) are formatted as Text ("@"), so writing a
# number straight into .Value would store it as a text string. Temporarily
# switch those cells to General, write the numeric 0, then restore the
# original Text format - matching how the cells are stored in the source
# file (numeric <v>0</v> even though the column is Text-formatted).
$lRange = $ws.Range("L459:L461")
$lRange.NumberFormat = "General"
$ws.Range("L459").Value = 0
$ws.Range("L460").Value = 0
$ws.Range("L461").Value = 0
$lRange.NumberFormat = "@"

$mRange = $ws.Range("M459:M461")
$mRange.NumberFormat = "General"
$ws.Range("M459").Value = 0
$ws.Range("M460").Value = 0
$ws.Range("M461").Value = 0
$mRange.NumberFormat = "@"
